$d = $word.ActiveDocument

# 1) Drop " la domiciliul solicitantuluisolicitantei" so the sentence reads
#    "...ocazia verificarilor si care va fi anexat..."
$d.Content.Find.Execute(
    "verificărilor la domiciliul solicitantuluisolicitantei și",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "verificărilor și", 2)

# 2) Replace "...politie care executa verificarile." with the finalized
#    wording that routes the report/proces-verbal onward.
$d.Content.Find.Execute(
    "poliție care execută verificările.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "poliție. Atât raportul cât și procesul-verbal vor fi înaintate către Serviciul Arme, Explozivi și Substanțe Periculoase.",
    2)
